$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 47, shifting existing rows 47-50 down to 48-51
$ws.Rows.Item(47).Insert()

# Fill in the new row 47 with the new weekly record
$ws.Range("A47").Value = 9
$ws.Range("B47").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C47").Value = "Metropolitana"
$ws.Range("D47").Value = 45021
$ws.Range("E47").Value = 13
$ws.Range("F47").Value = 100112010
$ws.Range("G47").Value = "Achicoria"
$ws.Range("H47").Value = "Sin especificar"
$ws.Range("I47").Value = "Primera"
$ws.Range("J47").Value = 50
$ws.Range("K47").Value = 7000
$ws.Range("L47").Value = 7000
$ws.Range("M47").Value = 7000
$ws.Range("N47").Value = "$/caja 16 unidades"
$ws.Range("O47").Value = "Provincia de Quillota"
$ws.Range("P47").Value = 438
$ws.Range("Q47").Value = 16
$ws.Range("R47").Value = "Hortaliza"
